$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Courses": add two new course rows (BITSF111, BITSF112)
# ---------------------------------------------------------------------
$courses = $wb.Worksheets.Item("Courses")

# Row 5 - BITSF111 / Thermodynamics
$courses.Cells.Item(5, 1).Value = "BITSF111"
$courses.Cells.Item(5, 2).Value = "Thermodynamics"
$courses.Cells.Item(5, 3).Value = 45271
$courses.Cells.Item(5, 3).NumberFormat = "yyyy-mm-dd h:mm:ss"
$courses.Cells.Item(5, 4).Value = "AN"

# Row 6 - BITSF112 / Technical Report Writting
$courses.Cells.Item(6, 1).Value = "BITSF112"
$courses.Cells.Item(6, 2).Value = "Technical Report Writting"
$courses.Cells.Item(6, 3).Value = 45269
$courses.Cells.Item(6, 3).NumberFormat = "yyyy-mm-dd h:mm:ss"
$courses.Cells.Item(6, 4).Value = "AN"

# ---------------------------------------------------------------------
# Sheet "Sections": add the matching section row for BITSF112
# ---------------------------------------------------------------------
$sections = $wb.Worksheets.Item("Sections")

# Row 11 already held MATHF111/T1 data but without explicit cell
# formatting; re-enter it so the row now carries the same formatting as
# the rest of the table.
$sections.Range("A11:F11").ClearContents()
$sections.Cells.Item(11, 1).Value = "MATHF111"
$sections.Cells.Item(11, 2).Value = "T1"
$sections.Cells.Item(11, 3).Value = 1
$sections.Cells.Item(11, 4).Value = "Th"
$sections.Cells.Item(11, 5).Value = 6168
$sections.Cells.Item(11, 6).Value = "Divyum"

# Row 12 - BITSF112, P1 section
$sections.Cells.Item(12, 1).Value = "BITSF112"
$sections.Cells.Item(12, 2).Value = "P1"
$sections.Cells.Item(12, 3).Value = 1
$sections.Cells.Item(12, 4).Value = "Th"
$sections.Cells.Item(12, 5).Value = 6164
$sections.Cells.Item(12, 6).Value = "Raj"
